$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.078.37'
$ws.Range("E2").Value = '  +6.07%  '
$ws.Range("D3").Value = '3.109.71'
$ws.Range("E3").Value = '  +3.89%  '
$ws.Range("E4").Value = '  -0.01%  '
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.72'
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = '  +3.78%  '
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.26'
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = '  +3.93%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.101.55'
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("E9").Value = '  +2.02%  '
$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = '  +13.31%  '
$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.77'
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = '  +8.61%  '
$ws.Range("E12").Value = '  +3.00%  '
$ws.Range("E13").Value = '  +7.60%  '
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.53'
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = '  +5.11%  '
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").Value = '3.623.38'
$ws.Range("E16").Value = '  +3.84%  '
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("D18").Value = '62.980.19'
$ws.Range("E18").Value = '  +5.90%  '
$ws.Range("D19").Value = '3.108.27'
$ws.Range("E19").Value = '  +3.86%  '
$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '466.05'
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = '  +7.06%  '
$ws.Range("E21").Value = '  +4.24%  '
$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.727'
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = '  +1.04%  '
$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = '  +7.07%  '
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("E26").Value = '  -0.12%  '
$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.37'
$ws.Range("D27").Style = $style_D27
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  +5.07%  '
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  +9.43%  '
$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.94'
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  +4.52%  '
$ws.Range("E33").Value = '  +2.56%  '
$ws.Range("D34").Value = '0.0₃0862'
$ws.Range("E34").Value = '  +10.78%  '
$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.42'
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = '  +16.29%  '
$ws.Range("E36").Value = '  +4.70%  '
$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.33'
$ws.Range("D37").Style = $style_D37
$ws.Range("E37").Value = '  +19.74%  '
$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.05'
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = '  +2.87%  '
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.81'
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = '  +4.05%  '
$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '432.75'
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = '  +7.73%  '
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.73'
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = '  +0.90%  '
$ws.Range("D42").Value = '2.923.12'
$ws.Range("E43").Value = '  +4.36%  '
$ws.Range("E44").Value = '  +11.73%  '
$ws.Range("E45").Value = '  +6.00%  '
$ws.Range("E46").Value = '  +7.46%  '
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.22'
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  +1.94%  '
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.77'
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("E50").Value = '  +0.76%  '
$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.51'
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = '  +4.51%  '
